$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- New header labels (row 1), columns AP:AW ---
$ws.Range("AP1").Value = "maxVolumePotSTb1"
$ws.Range("AQ1").Value = "maxVolumePotSTb2"
$ws.Range("AR1").Value = "maxVolumePotSTb3"
$ws.Range("AS1").Value = "maxVolumePotSTb4"
$ws.Range("AT1").Value = "minVolumePotSTb1"
$ws.Range("AU1").Value = "minVolumePotSTb2"
$ws.Range("AV1").Value = "minVolumePotSTb3"
$ws.Range("AW1").Value = "minVolumePotSTb4"

# --- New data values, rows 2:8, columns AP:AW ---
# row 2 - PHES Water reservoir (no explicit number format on AP:AS)
$ws.Range("AP2").Value = 2123749.4
$ws.Range("AQ2").Value = 1951920.27
$ws.Range("AR2").Value = 1354741.65
$ws.Range("AS2").Value = 1750162.93
$ws.Range("AT2:AW2").NumberFormat = "0.0"
$ws.Range("AT2").Value = 0
$ws.Range("AU2").Value = 0
$ws.Range("AV2").Value = 0
$ws.Range("AW2").Value = 0

# row 3 - Li-ion battery storage
$ws.Range("AP3:AW3").NumberFormat = "0.0"
$ws.Range("AP3").Value = 40000
$ws.Range("AQ3").Value = 40000
$ws.Range("AR3").Value = 40000
$ws.Range("AS3").Value = 40000
$ws.Range("AT3").Value = 0
$ws.Range("AU3").Value = 0
$ws.Range("AV3").Value = 0
$ws.Range("AW3").Value = 0

# row 4 - Hydrogen storage (high pressure)
$ws.Range("AP4:AW4").NumberFormat = "0.0"
$ws.Range("AP4").Value = 43200000
$ws.Range("AQ4").Value = 43200000
$ws.Range("AR4").Value = 43200000
$ws.Range("AS4").Value = 43200000
$ws.Range("AT4").Value = 0
$ws.Range("AU4").Value = 0
$ws.Range("AV4").Value = 0
$ws.Range("AW4").Value = 0

# row 5 - Methane storage tank
$ws.Range("AP5:AW5").NumberFormat = "0.0"
$ws.Range("AP5").Value = 50000
$ws.Range("AQ5").Value = 50000
$ws.Range("AR5").Value = 50000
$ws.Range("AS5").Value = 50000
$ws.Range("AT5").Value = 0
$ws.Range("AU5").Value = 0
$ws.Range("AV5").Value = 0
$ws.Range("AW5").Value = 0

# row 6 - Vanadium Redox Flow battery electrolyte
$ws.Range("AP6:AW6").NumberFormat = "0.0"
$ws.Range("AP6").Value = 50000
$ws.Range("AQ6").Value = 50000
$ws.Range("AR6").Value = 50000
$ws.Range("AS6").Value = 50000
$ws.Range("AT6").Value = 0
$ws.Range("AU6").Value = 0
$ws.Range("AV6").Value = 0
$ws.Range("AW6").Value = 0

# row 7 - Molten salts (high temperature heat)
$ws.Range("AP7:AW7").NumberFormat = "0.0"
$ws.Range("AP7").Value = 50000
$ws.Range("AQ7").Value = 50000
$ws.Range("AR7").Value = 50000
$ws.Range("AS7").Value = 50000
$ws.Range("AT7").Value = 0
$ws.Range("AU7").Value = 0
$ws.Range("AV7").Value = 0
$ws.Range("AW7").Value = 0

# row 8 - CO2 (captured) storage
$ws.Range("AP8:AW8").NumberFormat = "0.0"
$ws.Range("AP8").Value = 1000000000
$ws.Range("AQ8").Value = 1000000000
$ws.Range("AR8").Value = 1000000000
$ws.Range("AS8").Value = 1000000000
$ws.Range("AT8").Value = 0
$ws.Range("AU8").Value = 0
$ws.Range("AV8").Value = 0
$ws.Range("AW8").Value = 0

# --- View state: final selection is AU16 ---
$null = $ws.Range("AU16").Select()
